$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 33) holds a date value that was refreshed
# from 45620 (2024-11-24) to 45621 (2024-11-25).
for ($r = 2; $r -le 33; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45620) {
        $cell.Value2 = 45621
    }
}
